# Update "horarios 141" workbook with the 11:45 (08:45) scrape for LP1912, LP1912-215 and 6203-6173.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912   (cols: A blank | B Hora_Scrap | C Hora_Llegada | D Linea | E Minutos | F Parada | G Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 08:45:49"
$ws1.Range("A3").Value = "Total filas: 74"

$sheet1Rows = @(
    @("08:45:38", "08:51", "16_SANTA ANA",         6, "LP1912", "30/12/2025"),
    @("08:45:38", "08:53", "10_OLMOS",              8, "LP1912", "30/12/2025"),
    @("08:45:38", "09:01", "215A_EL PATO",         16, "LP1912", "30/12/2025"),
    @("08:45:38", "09:03", "11_ETCHEVERRY",        18, "LP1912", "30/12/2025"),
    @("08:45:38", "09:08", "23_HERNANDEZ",         23, "LP1912", "30/12/2025"),
    @("08:45:38", "09:10", "16_P MOR-SANTA ANA",   25, "LP1912", "30/12/2025"),
    @("08:45:38", "09:13", "10_OLMOS",             28, "LP1912", "30/12/2025"),
    @("08:45:38", "09:16", "27_EL RETIRO",         31, "LP1912", "30/12/2025"),
    @("08:45:38", "09:21", "26_HERNANDEZ",         36, "LP1912", "30/12/2025"),
    @("08:45:38", "09:22", "16_SANTA ANA",         37, "LP1912", "30/12/2025"),
    @("08:45:38", "09:23", "11_ETCHEVERRY",        38, "LP1912", "30/12/2025"),
    @("08:45:38", "09:32", "15_ABASTO",            47, "LP1912", "30/12/2025"),
    @("08:45:38", "09:33", "10_OLMOS",             48, "LP1912", "30/12/2025"),
    @("08:45:38", "09:42", "215C_EL PATO",         57, "LP1912", "30/12/2025"),
    @("08:45:38", "09:43", "14_ABASTO",            58, "LP1912", "30/12/2025"),
    @("08:45:38", "09:46", "23_HERNANDEZ",         61, "LP1912", "30/12/2025"),
    @("08:45:38", "09:52", "15_ABASTO",            67, "LP1912", "30/12/2025"),
    @("08:45:38", "10:03", "11_ETCHEVERRY",        78, "LP1912", "30/12/2025"),
    @("08:45:38", "10:10", "16_P MOR-SANTA ANA",   85, "LP1912", "30/12/2025"),
    @("08:45:38", "10:12", "15_ABASTO",            87, "LP1912", "30/12/2025"),
    @("08:45:38", "10:22", "17_ROMERO",            97, "LP1912", "30/12/2025")
)

$r = 55
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215   (cols: A blank | B Fecha | C Hora_Scrap | D Hora_Llegada | E Linea | F Minutos | G Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 08:45:49"
$ws2.Range("A3").Value = "Total filas: 10"

$sheet2Rows = @(
    @("30/12/2025", "08:45:38", "09:01", "215A_EL PATO", 16, "LP1912"),
    @("30/12/2025", "08:45:38", "09:42", "215C_EL PATO", 57, "LP1912")
)

$r = 10
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173   (cols: A blank | B Fecha | C Hora_Scrap | D Hora_Llegada | E Linea | F Minutos | G Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 08:45:49"
$ws3.Range("A3").Value = "Total filas: 11"

$sheet3Rows = @(
    @("30/12/2025", "08:45:49", "08:47", "215A_LA PLATA",           2, "L6173"),
    @("30/12/2025", "08:45:44", "09:09", "215D_LA PLATA",          24, "L6203"),
    @("30/12/2025", "08:45:49", "10:02", "215B_LP-P MOR-40 Y 115", 77, "L6173")
)

$r = 10
foreach ($row in $sheet3Rows) {
    $ws3.Cells.Item($r, 2).Value = $row[0]
    $ws3.Cells.Item($r, 3).Value = $row[1]
    $ws3.Cells.Item($r, 4).Value = $row[2]
    $ws3.Cells.Item($r, 5).Value = $row[3]
    $ws3.Cells.Item($r, 6).Value = $row[4]
    $ws3.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
